# Rewrite the "KEY ACHIEVEMENTS AND IMPACT" -> "Impact" bullet list so it
# matches the new, impact-focused accomplishment statements (6 bullets -> 4).

$d = $word.ActiveDocument

# Locate the "KEY ACHIEVEMENTS AND IMPACT" heading paragraph so the edits
# below are scoped to that section only (some bullet text is duplicated
# elsewhere in the resume, e.g. under Professional Experience).
$sectionIndex = 0
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*KEY ACHIEVEMENTS AND IMPACT*") {
        $sectionIndex = $i
        break
    }
}

if ($sectionIndex -eq 0) {
    throw "Could not find 'KEY ACHIEVEMENTS AND IMPACT' heading"
}

# Layout directly under the heading:
#   sectionIndex + 0 : "KEY ACHIEVEMENTS AND IMPACT" (heading)
#   sectionIndex + 1 : "Impact" (sub-heading)
#   sectionIndex + 2 : bullet 1 - race coding discovery
#   sectionIndex + 3 : bullet 2 - 87% prediction accuracy / polling margins
#   sectionIndex + 4 : bullet 3 - cloud data warehouse
#   sectionIndex + 5 : bullet 4 - redistricting platform (to remove)
#   sectionIndex + 6 : bullet 5 - longitudinal data analysis methods
#   sectionIndex + 7 : bullet 6 - ETL pipelines

$bullet1Index = $sectionIndex + 2
$bullet2Index = $sectionIndex + 3
$bullet3Index = $sectionIndex + 4
$bullet4Index = $sectionIndex + 5
$bullet5Index = $sectionIndex + 6
$bullet6Index = $sectionIndex + 7

# --- Sanity-check each paragraph before touching it ---
$b1 = $d.Paragraphs($bullet1Index).Range
if ($b1.Text -notlike "*Discovered systematic race coding errors*") {
    throw "Unexpected text at bullet 1"
}
$b2 = $d.Paragraphs($bullet2Index).Range
if ($b2.Text -notlike "*Achieved 87% prediction accuracy*") {
    throw "Unexpected text at bullet 2"
}
$b3 = $d.Paragraphs($bullet3Index).Range
if ($b3.Text -notlike "*Built cloud-based data warehouse*") {
    throw "Unexpected text at bullet 3"
}
$b4 = $d.Paragraphs($bullet4Index).Range
if ($b4.Text -notlike "*Built redistricting platform*") {
    throw "Unexpected text at bullet 4"
}
$b5 = $d.Paragraphs($bullet5Index).Range
if ($b5.Text -notlike "*Developed longitudinal data analysis methods*") {
    throw "Unexpected text at bullet 5"
}
$b6 = $d.Paragraphs($bullet6Index).Range
if ($b6.Text -notlike "*Designed ETL pipelines*") {
    throw "Unexpected text at bullet 6"
}

# --- Rewrite the three bullets that are kept in place ---
$b1.Text = [char]8226 + " Predictive excellence: Achieved 87% voter turnout accuracy vs. 71% industry standard"
$b2.Text = [char]8226 + " Reduced polling margins from " + [char]0xB1 + "4.2% to " + [char]0xB1 + "2.1%"
$b3.Text = [char]8226 + " Methodological advancement: Improved segmentation accuracy 34% and survey incidence 28%"

# --- Remove the two bullets that disappear entirely (redistricting platform,
#     longitudinal data analysis methods). Deleting twice at the same index
#     removes both, since subsequent paragraphs shift up after each delete. ---
$d.Paragraphs($bullet4Index).Range.Delete()
$d.Paragraphs($bullet4Index).Range.Delete()

# --- The final bullet (ETL pipelines) is now at $bullet4Index; rewrite it. ---
$bFinal = $d.Paragraphs($bullet4Index).Range
if ($bFinal.Text -notlike "*Designed ETL pipelines*") {
    throw "Unexpected text at final bullet position after deletion"
}
$bFinal.Text = [char]8226 + " Reduced polling costs while increasing quality"
